$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to stay text-typed so values like "0.999" or "6.71"
# are not auto-coerced to numbers by Excels input parser (matches original inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

# Price column (D) updates
$ws.Range('D2').Value = '63.741.99'
$ws.Range('D3').Value = '3.432.66'
$ws.Range('D5').Value = '577.28'
$ws.Range('D6').Value = '156.48'
$ws.Range('D7').Value = '0.999'
$ws.Range('D8').Value = '3.436.24'
$ws.Range('D11').Value = '0.122'
$ws.Range('D12').Value = '0.437'
$ws.Range('D13').Value = '4.022.47'
$ws.Range('D16').Value = '27.20'
$ws.Range('D17').Value = '63.821.78'
$ws.Range('D18').Value = '3.433.35'
$ws.Range('D19').Value = '6.42'
$ws.Range('D20').Value = '14.24'
$ws.Range('D21').Value = '8.49'
$ws.Range('D22').Value = '391.14'
$ws.Range('D23').Value = '0.999'
$ws.Range('D25').Value = '72.11'
$ws.Range('D27').Value = '9.50'
$ws.Range('D28').Value = '0.181'
$ws.Range('D30').Value = '6.71'
$ws.Range('D31').Value = '1.38'
$ws.Range('D32').Value = '2.04'
$ws.Range('D34').Value = '23.46'
$ws.Range('D36').Value = '6.79'
$ws.Range('D37').Value = '1.50'
$ws.Range('D38').Value = '158.89'
$ws.Range('D39').Value = '28.02'
$ws.Range('D40').Value = '0.0780'
$ws.Range('D41').Value = '1.87'
$ws.Range('D42').Value = '2.913.65'
$ws.Range('D43').Value = '0.0320'
$ws.Range('D44').Value = '0.769'
$ws.Range('D45').Value = '41.86'
$ws.Range('D46').Value = '4.36'
$ws.Range('D48').Value = '3.480.04'
$ws.Range('D49').Value = '22.56'
$ws.Range('D51').Value = '295.41'

# Volume / name / link column updates
$ws.Range('E2').Value = '  +5.98%  '
$ws.Range('E3').Value = '  +7.17%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +7.24%  '
$ws.Range('E6').Value = '  +7.35%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +7.04%  '
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').Value = '  +3.06%  '
$ws.Range('E11').Value = '  +8.47%  '
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E13').Value = '  +7.18%  '
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('E15').Value = '  +8.03%  '
$ws.Range('E16').Value = '  +5.41%  '
$ws.Range('E17').Value = '  +6.06%  '
$ws.Range('E18').Value = '  +6.72%  '
$ws.Range('E19').Value = '  +1.93%  '
$ws.Range('E20').Value = '  +7.45%  '
$ws.Range('E22').Value = '  +5.42%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +3.08%  '
$ws.Range('E26').Value = '  +22.21%  '
$ws.Range('E27').Value = '  +10.20%  '
$ws.Range('E28').Value = '  +6.66%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +9.33%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E31').Value = '  +16.76%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E32').Value = '  +7.74%  '
$ws.Range('E33').Value = '  +8.86%  '
$ws.Range('E34').Value = '  +4.43%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  +3.24%  '
$ws.Range('E37').Value = '  +8.99%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  +5.96%  '
$ws.Range('E40').Value = '  +10.10%  '
$ws.Range('E41').Value = '  +10.12%  '
$ws.Range('E42').Value = '  +4.29%  '
$ws.Range('E43').Value = '  +2.36%  '
$ws.Range('E44').Value = '  +6.91%  '
$ws.Range('E45').Value = '  +4.87%  '
$ws.Range('E46').Value = '  +3.44%  '
$ws.Range('E47').Value = '  +10.47%  '
$ws.Range('E48').Value = '  +7.26%  '
$ws.Range('E49').Value = '  +8.98%  '

# Reset style of touched D cells back to Normal (removes the temporary text numFmt)
# while keeping the stored values as text, so no stray "s" attribute is introduced.
$ws.Range("D2:D51").Style = "Normal"
